# Remove the <w:contextualSpacing w:val="0"/> element from every paragraph's
# paragraph properties (<w:pPr>) throughout the document body.
#
# The Word object model exposed by this runtime has no typed
# ParagraphFormat.ContextualSpacing property, so we perform the edit at the
# OOXML level: for each paragraph, pull its WordprocessingML via
# Range.WordOpenXML, strip the <w:contextualSpacing/> child (if present) from
# the captured <w:p>...</w:p> fragment, and write the trimmed markup back with
# Range.InsertXML (which replaces that range's contents in place, preserving
# paragraph count/order/text).

$d = $word.ActiveDocument

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
  $para = $d.Paragraphs($i)
  $openXml = $para.Range.WordOpenXML

  # The WordOpenXML payload wraps the whole package; pull out just this
  # paragraph's own <w:p>...</w:p> fragment (the first one - any synthetic
  # end-of-document paragraph that rides along on the final Paragraphs(Count)
  # range comes after it and is left untouched).
  if ($openXml -match '(?s)(<w:p\b.*?</w:p>)') {
    $pXml = $matches[1]

    if ($pXml -match '<w:contextualSpacing\b[^>]*/>') {
      $newXml = $pXml -replace '\s*<w:contextualSpacing\b[^>]*/>', ''
      $para.Range.InsertXML($newXml) | Out-Null
    }
  }
}

Write-Output "removed contextualSpacing from $count paragraphs"
